# Adjust labour type simulation
# - Map adjusted to new codes of model industries
# - Code adjusted to new structure of library code

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# New industry codes (row 1, columns B:AJ) replacing the old i001..i058 codes
$newCodes = @(
    "iPARI",
    "iWHEA",
    "iOCER",
    "iFVEG",
    "iOILS",
    "iSUGB",
    "iFIBR",
    "iOTHC",
    "iANIM",
    "iFORE",
    "iFISH",
    "iFOSM",
    "iOTHM",
    "iFBTO",
    "iTXWO",
    "iCOKE",
    "iREFN",
    "iCHEM",
    "iRUBP",
    "iNMMP",
    "iMETP",
    "iELEC",
    "iMACH",
    "iELCF",
    "iELCG",
    "iTRDI",
    "iHWAT",
    "iWATR",
    "iCONS",
    "iTRAD",
    "iHORE",
    "iTRAN",
    "iREBA",
    "iPUBO",
    "iWAST"
)

for ($i = 0; $i -lt $newCodes.Length; $i++) {
    # Column B is index 2
    $ws.Cells.Item(1, 2 + $i).Value = $newCodes[$i]
}

# Reset the frozen-pane view: scroll back to the top-left of the scrollable
# area and clear the stale selection on the bottom-right pane.
$ws.Activate()
$excel.ActiveWindow.ScrollRow = 2
$excel.ActiveWindow.ScrollColumn = 2
$ws.Range("B2").Select()
